$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handoff"
#
# The df2e5e7c-... entry is gone (handed off / completed), and the
# 6ae3d7f8-... entry's status flips from "Handed back: in sync with en-US"
# to "Ready for handoff" with refreshed handoff timestamps. Concretely,
# on each of the three worksheets (Overview, zh-cn, de-de):
#   - row 3 (the df2e5e7c-... record) is removed entirely, including its
#     hyperlinks
#   - row 2's status cell(s) change to "Ready for handoff"
#   - row 2's handoff-datetime cell is bumped to a new timestamp
# ---------------------------------------------------------------------------

function Remove-HyperlinksInRow {
    param($ws, [int]$row)
    $keepGoing = $true
    while ($keepGoing) {
        $keepGoing = $false
        foreach ($hl in $ws.Hyperlinks) {
            if ($hl.Range.Row -eq $row) {
                $hl.Delete()
                $keepGoing = $true
                break
            }
        }
    }
}

# --- Sheet "Overview" ------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-39-12 18:39:43"

Remove-HyperlinksInRow $wsOverview 3
$wsOverview.Rows.Item(3).Delete()

# --- Sheet "zh-cn" -----------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-12 18:39:40"

Remove-HyperlinksInRow $wsZhCn 3
$wsZhCn.Rows.Item(3).Delete()

# --- Sheet "de-de" -----------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-12 18:39:43"

Remove-HyperlinksInRow $wsDeDe 3
$wsDeDe.Rows.Item(3).Delete()
